$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 17
$ws.Range("D8").Value = 4
$ws.Range("D9").Value = 49
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("D19").Value = 1
$ws.Range("D20").Value = 15
$ws.Range("D21").Value = 1
$ws.Range("D22").Value = 1
